# Applies the diff:
#   - removes the stray "_GoBack" bookmark that used to sit right after the
#     second diagram image
#   - reworks "Implementar la gestión de artículos comprados." into
#     "Implementar la gestión de artículos por comprar." across several
#     runs, with the (moved) "_GoBack" bookmark now sitting right before the
#     final period, matching how Word leaves its "last edit" bookmark behind
#     after an in-place text edit.

$d = $word.ActiveDocument

# Locate the target paragraph ("Implementar la gestión de artículos
# comprados.") by content instead of a hard-coded index.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Implementar la gestión de artículos comprados.*") {
        $target = $p
    }
}

if ($target -ne $null) {
    # --- Phase 1: plain content edits (kept to as few distinct run-touching
    # operations as possible, since replacing/deleting text tends to merge
    # adjacent same-formatted runs back together). ---

    # "comprados" -> "comprar"
    $rFix = $target.Range.Duplicate
    $rFix.Find.Execute("comprados", $true, $false, $false, $false, $false, $true, 1, $false, "comprar", 2)

    # insert "por " right after "artículos "
    $rArt = $target.Range.Duplicate
    $rArt.Find.Execute("artículos ")
    $porStart = $rArt.End
    $insPor = $d.Range($porStart, $porStart)
    $insPor.InsertAfter("por ")

    # --- Phase 2: split the paragraph into the final run layout. Wrapping a
    # (non-empty) temporary bookmark around a substring forces Word to end
    # the surrounding runs at the bookmark's edges; deleting the bookmark
    # right afterwards leaves the runs split apart without merging them
    # back together, as long as no further text is inserted/deleted. ---

    # split out "por " into its own run
    $porRange = $d.Range($porStart, $porStart + 4)
    $d.Bookmarks.Add("_TmpSplitPor", $porRange)
    $d.Bookmarks("_TmpSplitPor").Delete()

    # split "compra" away from the trailing "r."
    $rCompra = $target.Range.Duplicate
    $rCompra.Find.Execute("compra")
    $d.Bookmarks.Add("_TmpSplitCompra", $rCompra)
    $d.Bookmarks("_TmpSplitCompra").Delete()

    # finally, drop the real "_GoBack" bookmark right before the trailing
    # period -- this both splits "r" from "." and (because a document can
    # only have one "_GoBack" bookmark) removes the old one that used to
    # live after the second diagram's picture.
    $rEnd = $target.Range.Duplicate
    $rEnd.Find.Execute("comprar.")
    $dotPos = $rEnd.End - 1
    $insDot = $d.Range($dotPos, $dotPos)
    $d.Bookmarks.Add("_GoBack", $insDot)
}
